# Generate Report for Handback
#
# Refresh the "Latest HO Xliff Generate Date" on the Overview sheet and the
# corresponding handoff/handback timestamps on the per-locale (zh-cn, de-de)
# sheets for the 5ad28dad-3920-4125-bac1-dd1defa03890.md row, reflecting a
# newer handback report generation run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
# Row 3 corresponds to 5ad28dad-3920-4125-bac1-dd1defa03890.md
$overview.Range("G3").Value = "2016-08-31 15:02:23"

# --- zh-cn sheet ------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
# Row 3 corresponds to 5ad28dad-3920-4125-bac1-dd1defa03890.md
$zhcn.Range("H3").Value = "2016-08-31 15:02:10"
$zhcn.Range("K3").Value = "2016-08-31 15:02:51"

# --- de-de sheet ------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
# Row 3 corresponds to 5ad28dad-3920-4125-bac1-dd1defa03890.md
$dede.Range("H3").Value = "2016-08-31 15:02:23"
$dede.Range("K3").Value = "2016-08-31 15:03:00"
